$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (current Match-ref is C, Primers is D)
# This will shift Primers..primer-mismatch-rate one column to the right (D->E ... H->I)
$ws.Columns.Item(4).Insert()

# Set the new column D width to (approximately) match column C's width
$ws.Columns.Item(4).ColumnWidth = 8.33

# Header for new column D: "Segmented"
$ws.Range("D1").Value = "Segmented"

# Fill in boolean "Segmented" values for rows 2-5
$ws.Range("D2").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $true
$ws.Range("D5").Value = $false

# Update "Match-ref" column (C) values per new data
$ws.Range("C2").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = $true
$ws.Range("C5").Value = $true

# Update selection to D3
$ws.Range("D3").Select()
